$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price cells to remain text (preserves original text formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "25.944.39"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "1.638.73"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "214.72"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "0.06360"
$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  +1.60%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").Value = "1.638.45"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").Value = "0.5471"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "0.0₅7751"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("D16").Value = "64.22"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "25.961.20"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").Value = "4.463"

$ws.Range("D20").Value = "195.92"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "9.956"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "6.137"
$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.56%  "

$ws.Range("D24").Value = "1.896"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").Value = "143.04"
$ws.Range("E25").Value = "  +0.88%  "

$ws.Range("D26").Value = "0.1258"
$ws.Range("E26").Value = "  +10.12%  "

$ws.Range("D27").Value = "6.852"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("D28").Value = "15.62"
$ws.Range("E28").Value = "  -0.72%  "

$ws.Range("D29").Value = "1.237"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").Value = "0.04888"
$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("D31").Value = "3.253"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").Value = "3.205"
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("D33").Value = "1.556"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("D34").Value = "2.375"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").Value = "0.9164"
$ws.Range("E35").Value = "  +2.34%  "

$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").Value = "1.136.80"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "0.5528"
$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("D41").Value = "5.593"
$ws.Range("E41").Value = "  -0.85%  "

$ws.Range("D42").Value = "0.8044"
$ws.Range("E42").Value = "  -1.55%  "

$ws.Range("D43").Value = "98.50"
$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("E44").Value = "  -9.84%  "

$ws.Range("D45").Value = "1.771.25"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("D46").Value = "0.4498"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").Value = "0.05182"
$ws.Range("E49").Value = "  +1.97%  "

$ws.Range("D50").Value = "7.491"
$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("E51").Value = "  -0.45%  "
